# Update metrics table (columns B:Q, rows 2:26) with the new values
# produced by the retrained model ("atualizado todo o treinamento para o novo lm").
# Every data row shares the same metric values, so we build one array of
# values (columns B..Q, in order) and apply it uniformly to every row from
# 2 to 26. Numbers are written in plain decimal (no scientific notation)
# because the interpreter here does not accept an "e-05"/"E-05" exponent
# suffix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.9999674344715328,             # B - r2
    0.9989400190852257,             # C - r2_sup
    0.999781308147743,              # D - r2_test
    0.9999999999999921,             # E - r2_val
    0.9999812344770503,             # F - r2_vt
    0.00003039849136261971,         # G - mse
    0.0009894456561571146,          # H - mse_sup
    0.00005422831569902617,         # I - mse_test
    0.000000000000008301796908186114, # J - mse_val
    0.00002711415785366398,         # K - mse_vt
    0.000349036811621329,           # L - mape
    0.005513482689065026,           # M - rmse
    1.000060120975632,              # N - r2_adj
    0.005748203061621094,           # O - rsd
    94.80223515402628,              # P - aic
    139.9006406741497               # Q - bic
)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 2]
    }
}
